# Apply "revisi lagi P3P Agustus 2020" edit: extend calendar table with weeks 12-20 (rows 8-42, columns F-I)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F8").Value = 12
$ws.Range("G8").Value = "Rabu"
$ws.Range("H8").Value = "Ams 12"
$ws.Range("I8").Value = "Mzm 75"

$ws.Range("H9").Value = "Luk 3"
$ws.Range("I9").Value = "Hak 13"

$ws.Range("H10").Value = "1 Kor 10"
$ws.Range("I10").Value = "Yeh 6"

$ws.Range("F12").Value = 13
$ws.Range("G12").Value = "Kamis"
$ws.Range("H12").Value = "Ams 13"
$ws.Range("I12").Value = "Mzm 76"

$ws.Range("H13").Value = "Luk 4"
$ws.Range("I13").Value = "Hak 14"

$ws.Range("H14").Value = "1 Kor 11"
$ws.Range("I14").Value = "Yeh 7"

$ws.Range("F16").Value = 14
$ws.Range("G16").Value = "Jumat"
$ws.Range("H16").Value = "Ams 14"
$ws.Range("I16").Value = "Mzm 77"

$ws.Range("H17").Value = "Luk 5"
$ws.Range("I17").Value = "Hak 15"

$ws.Range("H18").Value = "1 Kor 12"
$ws.Range("I18").Value = "Yeh 8"

$ws.Range("F20").Value = 15
$ws.Range("G20").Value = "Sabtu"
$ws.Range("H20").Value = "Ams 15"
$ws.Range("I20").Value = "Mzm 78"

$ws.Range("H21").Value = "Luk 6"
$ws.Range("I21").Value = "Hak 16"

$ws.Range("H22").Value = "1 Kor 13"
$ws.Range("I22").Value = "Yeh 9"

$ws.Range("F24").Value = 16
$ws.Range("G24").Value = "Minggu"
$ws.Range("H24").Value = "Ams 16"
$ws.Range("I24").Value = "Mzm 79"

$ws.Range("H25").Value = "Luk 7"
$ws.Range("I25").Value = "Hak 17"

$ws.Range("H26").Value = "1 Kor 14"
$ws.Range("I26").Value = "Yeh 10"

$ws.Range("F28").Value = 17
$ws.Range("G28").Value = "Senin"
$ws.Range("H28").Value = "Ams 17"
$ws.Range("I28").Value = "Mzm 80"

$ws.Range("H29").Value = "Luk 8"
$ws.Range("I29").Value = "Hak 18"

$ws.Range("H30").Value = "1 Kor 15"
$ws.Range("I30").Value = "Yeh 11"

$ws.Range("F32").Value = 18
$ws.Range("G32").Value = "Selasa"
$ws.Range("H32").Value = "Ams 18"
$ws.Range("I32").Value = "Mzm 81"

$ws.Range("H33").Value = "Luk 9"
$ws.Range("I33").Value = "Hak 19"

$ws.Range("H34").Value = "1 Kor 16"
$ws.Range("I34").Value = "Yeh 12"

$ws.Range("F36").Value = 19
$ws.Range("G36").Value = "Rabu"
$ws.Range("H36").Value = "Ams 19"
$ws.Range("I36").Value = "Mzm 82"

$ws.Range("H37").Value = "Luk 10"
$ws.Range("I37").Value = "Hak 20"

$ws.Range("I38").Value = "Yeh 13"
$ws.Range("H38").Value = "2 Kor 1"

$ws.Range("F40").Value = 20
$ws.Range("G40").Value = "Kamis "
$ws.Range("H40").Value = "Ams 20"
$ws.Range("I40").Value = "Mzm 83"

$ws.Range("H41").Value = "Luk 11"
$ws.Range("I41").Value = "Hak 21"

$ws.Range("H42").Value = "2 Kor 2"
$ws.Range("I42").Value = "Yeh 14"

# Update the selected cell to I43, matching the saved selection state in the workbook
$ws.Range("I43").Select()
